# Swap the two worker/account rows (16 and 17): the "EC" (estado de cuenta)
# for EMILER DEL CARMEN CASTRO BALLESTA and ARMANDO POMARES GUZMAN trade
# places - Armando's record now appears first (row 16), Emiler's second
# (row 17). Columns C (N° Doc), D (Nombre), E (Periodo Mora), F (Valor Mora)
# and G (Salario Basico) all move together; column B ("CC") is identical on
# both rows so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 16
$row2 = 17
$cols = @(3, 4, 5, 6, 7)

foreach ($c in $cols) {
    $cell1 = $ws.Cells.Item($row1, $c)
    $cell2 = $ws.Cells.Item($row2, $c)
    $tmp = $cell1.Value2
    $cell1.Value = $cell2.Value2
    $cell2.Value = $tmp
}
